# Cacheo de objetos desactivado en contenido
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear old rows 11 and 12 entirely, since data moves around / shrinks
$ws.Range("C11:E12").Clear()

# Insert new shared strings in the same order as the target workbook so the
# sharedStrings table comes out in the expected sequence:
#   TestMapeoAVortex, Hola manola, (Obj/ms), Obj->Vortex (obj/ms), ClaseParaProbarMapeo
$ws.Range("D15").Value = "TestMapeoAVortex"
$ws.Range("C17").Value = "Hola manola"
$ws.Range("D16").Value = "(Obj/ms)"
$ws.Range("D6").Value = "Obj->Vortex (obj/ms)"
$ws.Range("C18").Value = "ClaseParaProbarMapeo"

# Rest of the new / moved content
$ws.Range("E6").Value = "Vortex->Obj"

$ws.Range("C7").Value = "String"
$ws.Range("D7").Value = 252.91470000000001
$ws.Range("E7").Value = 252.91460000000001

$ws.Range("C8").Value = "Number"
$ws.Range("D8").Value = 263.20710000000003
$ws.Range("E8").Value = 263.20699999999999

$ws.Range("C9").Value = "List<String>"
$ws.Range("D9").Value = 236.62889999999999
$ws.Range("E9").Value = 236.62880000000001

$ws.Range("C10").Value = "modelo"
$ws.Range("D10").Value = 31.4053
$ws.Range("E10").Value = 31.4053

$ws.Range("C15").Value = "Test"

$ws.Range("D17").Value = 228.83295194508

$ws.Range("D18").Value = 1.0406585287169701

# Update column widths (runtime stores width = ColumnWidth + 0.8333333333333334,
# so subtract that offset to land on the target stored widths of 12 and 22)
$ws.Columns.Item(3).ColumnWidth = 11.166666666666666
$ws.Columns.Item(4).ColumnWidth = 21.166666666666668

# Update selection to match new active cell
$ws.Range("D11").Select()
